# Insert a new weekly price-report row for "Femacal de La Calera" (Haba)
# above the existing row 49, shifting all following rows down by one
# (dimension grows from A1:R76 to A1:R77), then populate the new row
# with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 49; everything below (old rows 49-76)
# shifts down to become rows 50-77.
$ws.Rows.Item(49).Insert()

# Fill in the newly inserted row 49 with the new weekly record.
$ws.Range("A49").Value = 3
$ws.Range("B49").Value = "Femacal de La Calera"
$ws.Range("C49").Value = "Coquimbo"
$ws.Range("D49").Value = 44488
$ws.Range("E49").Value = 5
$ws.Range("F49").Value = 100112026
$ws.Range("G49").Value = "Haba"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 80
$ws.Range("K49").Value = 9500
$ws.Range("L49").Value = 10000
$ws.Range("M49").Value = 9750
$ws.Range("N49").Value = "$/malla 25 kilos"
$ws.Range("O49").Value = "Provincia de Limarí"
$ws.Range("P49").Value = 390
$ws.Range("Q49").Value = 25
$ws.Range("R49").Value = "Hortaliza"
